# "March 24 update 3"
#
# Adds three new trailing columns (M: renewd, N: PlanID, O: iteration) to
# the table on Sheet1, matching the header style already used by columns
# B:L, and fills the new columns for the six data rows (2-7) with:
#   M = "before"
#   N = 20140231
#   O = 5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the bold / centered / bordered header style from the existing last
# header cell (L1) onto the three new header cells so they match the rest
# of the header row.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Populate the new columns for each data row.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"      # column M
    $ws.Cells.Item($r, 14).Value = 20140231      # column N
    $ws.Cells.Item($r, 15).Value = 5             # column O
}
